$wb = $excel.ActiveWorkbook

# This script updates the per-row market/profit figures (columns H-N) across
# several Leve "Profits" worksheets (ALC, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed market-board pricing data pulled by the scheduled runner.

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 237.44444
$ws.Range("I5").Value = 89.8
$ws.Range("J5").Value = 422
$ws.Range("K5").Value = 89.8
$ws.Range("L5").Value = 422
$ws.Range("M5").Value = 25.2
$ws.Range("N5").Value = -652
$ws.Range("H80").Value = 2613.4546
$ws.Range("J80").Value = 4239.6
$ws.Range("L80").Value = 12718.8
$ws.Range("N80").Value = -14714.8
$ws.Range("H83").Value = 2613.4546
$ws.Range("J83").Value = 4239.6
$ws.Range("L83").Value = 38156.4
$ws.Range("N83").Value = -48140.4
$ws.Range("H88").Value = 1719
$ws.Range("I88").Value = 697.75
$ws.Range("J88").Value = 1923.25
$ws.Range("K88").Value = 697.75
$ws.Range("L88").Value = 1923.25
$ws.Range("M88").Value = -291.75
$ws.Range("N88").Value = -2735.25
$ws.Range("H91").Value = 1719
$ws.Range("I91").Value = 697.75
$ws.Range("J91").Value = 1923.25
$ws.Range("K91").Value = 697.75
$ws.Range("L91").Value = 1923.25
$ws.Range("M91").Value = 706.25
$ws.Range("N91").Value = -4731.25
$ws.Range("H135").Value = 3531.6667
$ws.Range("I135").Value = 2995
$ws.Range("K135").Value = 26955
$ws.Range("M135").Value = -24420
$ws.Range("H137").Value = 2128.8064
$ws.Range("I137").Value = 1293.6471
$ws.Range("K137").Value = 3880.9413
$ws.Range("M137").Value = -1330.9413

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 480.66666
$ws.Range("I80").Value = 470
$ws.Range("J80").Value = 489.2
$ws.Range("K80").Value = 470
$ws.Range("L80").Value = 489.2
$ws.Range("M80").Value = 528
$ws.Range("N80").Value = -2485.2
$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 60000
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -9617
$ws.Range("N82").Value = -60766
$ws.Range("H83").Value = 480.66666
$ws.Range("I83").Value = 470
$ws.Range("J83").Value = 489.2
$ws.Range("K83").Value = 2350
$ws.Range("L83").Value = 2446
$ws.Range("M83").Value = 2642
$ws.Range("N83").Value = -12430
$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 60000
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -8674
$ws.Range("N85").Value = -62652
$ws.Range("H86").Value = 7833.3335
$ws.Range("I86").Value = 7833.3335
$ws.Range("K86").Value = 7833.3335
$ws.Range("M86").Value = -6710.3335
$ws.Range("H89").Value = 7833.3335
$ws.Range("I89").Value = 7833.3335
$ws.Range("K89").Value = 39166.6675
$ws.Range("M89").Value = -33550.6675
$ws.Range("H105").Value = 5495.6665
$ws.Range("I105").Value = 4737.3335
$ws.Range("K105").Value = 4737.3335
$ws.Range("M105").Value = -2990.3335

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5095.174
$ws.Range("I58").Value = 3900.6667
$ws.Range("J58").Value = 7334.875
$ws.Range("K58").Value = 3900.6667
$ws.Range("L58").Value = 7334.875
$ws.Range("M58").Value = -3697.6667
$ws.Range("N58").Value = -7740.875
$ws.Range("H86").Value = 11621161
$ws.Range("I86").Value = 17429242
$ws.Range("K86").Value = 17429242
$ws.Range("M86").Value = -17428119
$ws.Range("H89").Value = 11621161
$ws.Range("I89").Value = 17429242
$ws.Range("K89").Value = 87146210
$ws.Range("M89").Value = -87140594
$ws.Range("H94").Value = 3095.0625
$ws.Range("I94").Value = 2991.9092
$ws.Range("J94").Value = 3322
$ws.Range("K94").Value = 2991.9092
$ws.Range("L94").Value = 3322
$ws.Range("M94").Value = -2540.9092
$ws.Range("N94").Value = -4224
$ws.Range("H99").Value = 5653.5454
$ws.Range("I99").Value = 4798.0557
$ws.Range("K99").Value = 4798.0557
$ws.Range("M99").Value = -3300.0557
$ws.Range("H126").Value = 5653.5454
$ws.Range("I126").Value = 4798.0557
$ws.Range("K126").Value = 14394.1671
$ws.Range("M126").Value = -11924.1671
$ws.Range("H134").Value = 1665.8182
$ws.Range("J134").Value = 1749.75
$ws.Range("L134").Value = 5249.25
$ws.Range("N134").Value = -10319.25
$ws.Range("H136").Value = 5095.174
$ws.Range("I136").Value = 3900.6667
$ws.Range("J136").Value = 7334.875
$ws.Range("K136").Value = 11702.0001
$ws.Range("L136").Value = 22004.625
$ws.Range("M136").Value = -9152.000100000001
$ws.Range("N136").Value = -27104.625

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 490.5
$ws.Range("I6").Value = 546.4286
$ws.Range("J6").Value = 99
$ws.Range("K6").Value = 1639.2858
$ws.Range("L6").Value = 297
$ws.Range("M6").Value = -1526.2858
$ws.Range("N6").Value = -523
$ws.Range("H63").Value = 156
$ws.Range("J63").Value = 200
$ws.Range("L63").Value = 600
$ws.Range("N63").Value = -2098
$ws.Range("H66").Value = 156
$ws.Range("J66").Value = 200
$ws.Range("L66").Value = 1800
$ws.Range("N66").Value = -9288
$ws.Range("H87").Value = 969.5
$ws.Range("I87").Value = 969.5
$ws.Range("K87").Value = 2908.5
$ws.Range("M87").Value = -1660.5
$ws.Range("H90").Value = 969.5
$ws.Range("I90").Value = 969.5
$ws.Range("K90").Value = 8725.5
$ws.Range("M90").Value = -2485.5
$ws.Range("H95").Value = 14000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 14000
$ws.Range("K95").Value = 0
$ws.Range("M95").Value = 42000
$ws.Range("N95").Value = -46118
$ws.Range("H97").Value = 599.5
$ws.Range("I97").Value = 750
$ws.Range("J97").Value = 449
$ws.Range("K97").Value = 2250
$ws.Range("L97").Value = 1347
$ws.Range("M97").Value = -1754
$ws.Range("N97").Value = -2339
$ws.Range("H98").Value = 562.6667
$ws.Range("I98").Value = 525.6
$ws.Range("K98").Value = 1576.8
$ws.Range("M98").Value = -78.80000000000018
$ws.Range("H113").Value = 290
$ws.Range("I113").Value = 275
$ws.Range("J113").Value = 297.5
$ws.Range("K113").Value = 825
$ws.Range("L113").Value = 892.5
$ws.Range("M113").Value = 1345
$ws.Range("N113").Value = -5232.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 29999
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 29999
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = 29999
$ws.Range("N94").Value = -31351

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1250
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1250
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9999.5
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 9999.5
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 475.27274
$ws.Range("I113").Value = 490.8
$ws.Range("J113").Value = 320
$ws.Range("K113").Value = 1472.4
$ws.Range("L113").Value = 960
$ws.Range("M113").Value = 697.5999999999999
$ws.Range("N113").Value = -5300
$ws.Range("H132").Value = 44811.434
$ws.Range("J132").Value = 1997.5
$ws.Range("L132").Value = 5992.5
$ws.Range("N132").Value = -11052.5
$ws.Range("H141").Value = 93999.5
$ws.Range("J141").Value = 93999.5
$ws.Range("L141").Value = 93999.5
$ws.Range("N141").Value = -104359.5

Write-Output "Applied profit/market data updates to ALC, BSM, CRP, CUL, GSM, LTW, WVR"
